$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (lsh_text_out_categories)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "lsh_sheet_names"

# Set the width of column A (closest achievable snap to 24.33203125)
$newSheet.Columns.Item(1).ColumnWidth = 23.5

# Fill in the sheet-name data values first (so the header string is added to the
# shared-string table last, matching the order new unique strings were recorded)
$values = @(
    "Einangrun - af skjáborði Heilsu",
    "Lokaviðtal-Spurning úr forms",
    "Einstaklingar",
    "Komur og innlagnir",
    "Áhættuflokkur ofl úr hóp",
    "Fyrsta viðtal úr forms",
    "Spurningar úr forms Pivot",
    "NEWS score ",
    "EG_Skoða gagnagöt",
    "Eldra_Spurningar úr forms"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $values[$i]
}

# Header goes in last
$newSheet.Cells.Item(1, 1).Value = "sheet_name_raw"

# Match the saved selection/active cell on the new sheet
$newSheet.Range("F10").Select()

# Make the new sheet the active tab
$newSheet.Activate()
